$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.831993222236633
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 2.614305257797241
$ws.Range("D1").Value = 1.102808713912964
$ws.Range("E1").Value = 0.7519119381904602
